$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.119.62'
$ws.Range("E2").Value = '  -2.30%  '

$ws.Range("D3").Value = '1.898.82'
$ws.Range("E3").Value = '  -2.65%  '

$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").Value = '332.45'
$ws.Range("E5").Value = '  -2.92%  '

$ws.Range("E6").Value = '  -0.02%  '

$ws.Range("D7").Value = '0.4597'
$ws.Range("E7").Value = '  -3.58%  '

$ws.Range("D8").Value = '0.4127'
$ws.Range("E8").Value = '  -0.35%  '

$ws.Range("E9").Value = '  -0.01%  '

$ws.Range("D10").Value = '0.08008'
$ws.Range("E10").Value = '  -2.83%  '

$ws.Range("E11").Value = '  -2.37%  '

$ws.Range("D12").Value = '22.19'
$ws.Range("E12").Value = '  -2.24%  '

$ws.Range("D13").Value = '1.896.53'
$ws.Range("E13").Value = '  -2.36%  '

$ws.Range("D14").Value = '5.941'
$ws.Range("E14").Value = '  -3.90%  '

$ws.Range("D15").Value = '7.114'
$ws.Range("E15").Value = '  -4.02%  '

$ws.Range("E16").Value = '  +0.00%  '

$ws.Range("D17").Value = '89.13'
$ws.Range("E17").Value = '  -3.30%  '

$ws.Range("E18").Value = '  -2.95%  '

$ws.Range("D19").Value = '0.06566'
$ws.Range("E19").Value = '  -1.94%  '

$ws.Range("D20").Value = '17.62'
$ws.Range("E20").Value = '  -2.44%  '

$ws.Range("D21").Value = '1.000'
$ws.Range("E21").Value = '  -0.04%  '

$ws.Range("D22").Value = '29.091.73'
$ws.Range("E22").Value = '  -2.34%  '

$ws.Range("D23").Value = '5.494'
$ws.Range("E23").Value = '  -1.55%  '

$ws.Range("E24").Value = '  +1.24%  '

$ws.Range("D25").Value = '2.197'
$ws.Range("E25").Value = '  -2.83%  '

$ws.Range("D26").Value = '2.117.00'
$ws.Range("E26").Value = '  -2.84%  '

$ws.Range("D27").Value = '156.59'
$ws.Range("E27").Value = '  -3.19%  '

$ws.Range("E28").Value = '  -2.22%  '

$ws.Range("D29").Value = '2.120'
$ws.Range("E29").Value = '  -2.84%  '

$ws.Range("D30").Value = '5.616'
$ws.Range("E30").Value = '  -1.40%  '

$ws.Range("D31").Value = '116.94'

$ws.Range("D32").Value = '1.054'
$ws.Range("E32").Value = '  +4.27%  '

$ws.Range("D33").Value = '0.09393'
$ws.Range("E33").Value = '  -2.34%  '

$ws.Range("E34").Value = '  -3.88%  '

$ws.Range("D35").Value = '3.546'
$ws.Range("E35").Value = '  -3.85%  '

$ws.Range("D36").Value = '5.354'
$ws.Range("E36").Value = '  -3.00%  '

$ws.Range("D37").Value = '0.06089'
$ws.Range("E37").Value = '  -3.50%  '

$ws.Range("D38").Value = '0.02239'
$ws.Range("E38").Value = '  -3.32%  '

$ws.Range("D39").Value = '8.429'
$ws.Range("E39").Value = '  -0.80%  '

$ws.Range("D40").Value = '1.178'
$ws.Range("E40").Value = '  -0.64%  '

$ws.Range("D41").Value = '0.5838'
$ws.Range("E41").Value = '  -4.30%  '

$ws.Range("E42").Value = '  +0.04%  '

$ws.Range("E43").Value = '  -3.43%  '

$ws.Range("D44").Value = '10.12'
$ws.Range("E44").Value = '  -5.55%  '

$ws.Range("D45").Value = '2.356'
$ws.Range("E45").Value = '  -0.88%  '

$ws.Range("D46").Value = '1.250'
$ws.Range("E46").Value = '  -0.39%  '

$ws.Range("D47").Value = '0.07502'
$ws.Range("E47").Value = '  +2.20%  '

$ws.Range("D48").Value = '12.16'
$ws.Range("E48").Value = '  -2.47%  '

$ws.Range("D49").Value = '0.5544'
$ws.Range("E49").Value = '  -3.00%  '

$ws.Range("D50").Value = '1.924'
$ws.Range("E50").Value = '  -3.28%  '

$ws.Range("D51").Value = '112.34'
$ws.Range("E51").Value = '  -1.12%  '
